{"js": "// The document contains a numbered list of \"PERSON\" placeholder entries,\n// each shaped like: \"[[PERSON_N]] \u2013 <something>, <prep> <something>\".\n// Most paragraphs already use \"[[PERSON_N]]\" for the repeated <something>,\n// but a handful of paragraphs still spell out the literal (anonymized)\n// name instead of reusing the \"[[PERSON_N]]\" placeholder. This script\n// replaces those literal names with the matching \"[[PERSON_N]]\" token so\n// every paragraph consistently reuses the placeholder.\n\nconst replacements = [\n  { name: \"Ann\u011b Kubi\u0161tov\u00e9\", token: \"[[PERSON_12]]\" },\n  { name: \"Jan\u011b Lipov\u00e9\", token: \"[[PERSON_13]]\" },\n  { name: \"Ivan\u011b Hol\u00ednkov\u00e9\", token: \"[[PERSON_16]]\" },\n  { name: \"Elen\u011b Krbcov\u00e9\", token: \"[[PERSON_24]]\" },\n  { name: \"Dian\u011b Kaprov\u00e9\", token: \"[[PERSON_41]]\" },\n  { name: \"Rajn\u011b Divi\u0161ov\u00e9\", token: \"[[PERSON_47]]\" },\n  { name: \"Han\u011b Pivo\u0148kov\u00e9\", token: \"[[PERSON_49]]\" },\n  { name: \"Em\u011b Hru\u0161kov\u00e9\", token: \"[[PERSON_50]]\" },\n];\n\nconst body = context.document.body;\nconst searchResultsList = [];\n\nfor (const { name } of replacements) {\n  const results = body.search(name, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  searchResultsList.push(results);\n}\n\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const { token } = replacements[i];\n  const results = searchResultsList[i];\n  for (let j = 0; j < results.items.length; j++) {\n    results.items[j].insertText(token, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains a numbered list of \"PERSON\" placeholder entries,\n# each shaped like: \"[[PERSON_N]] - <something>, <prep> <something>\".\n# Most paragraphs already use \"[[PERSON_N]]\" for the repeated <something>,\n# but a handful of paragraphs still spell out the literal (anonymized)\n# name instead of reusing the \"[[PERSON_N]]\" placeholder. This script\n# replaces those literal names with the matching \"[[PERSON_N]]\" token so\n# every paragraph consistently reuses the placeholder.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Name = \"Ann\u011b Kubi\u0161tov\u00e9\";   Token = \"[[PERSON_12]]\" },\n    @{ Name = \"Jan\u011b Lipov\u00e9\";      Token = \"[[PERSON_13]]\" },\n    @{ Name = \"Ivan\u011b Hol\u00ednkov\u00e9\";  Token = \"[[PERSON_16]]\" },\n    @{ Name = \"Elen\u011b Krbcov\u00e9\";    Token = \"[[PERSON_24]]\" },\n    @{ Name = \"Dian\u011b Kaprov\u00e9\";    Token = \"[[PERSON_41]]\" },\n    @{ Name = \"Rajn\u011b Divi\u0161ov\u00e9\";   Token = \"[[PERSON_47]]\" },\n    @{ Name = \"Han\u011b Pivo\u0148kov\u00e9\";   Token = \"[[PERSON_49]]\" },\n    @{ Name = \"Em\u011b Hru\u0161kov\u00e9\";     Token = \"[[PERSON_50]]\" }\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Name\n    $find.Replacement.Text = $r.Token\n    $find.Forward = $true\n    $find.Wrap = 1\n    # wdFindContinue=1, wdReplaceAll=2\n    $find.Execute($r.Name, $false, $false, $false, $false, $false, $true, 1, $false, $r.Token, 2)\n}\n\nWrite-Output \"done\"\n"}
